$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.078.20"
$ws.Range("E2").Value = "'  +1.92%  "
$ws.Range("D3").Value = "'3.456.78"
$ws.Range("E3").Value = "'  +1.28%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("D5").Value = "'580.13"
$ws.Range("E5").Value = "'  +0.41%  "
$ws.Range("D6").Value = "'147.42"
$ws.Range("E6").Value = "'  +2.24%  "
$ws.Range("D7").Value = "'3.455.49"
$ws.Range("E7").Value = "'  +1.22%  "
$ws.Range("E8").Value = "'  -0.03%  "
$ws.Range("E9").Value = "'  +0.83%  "
$ws.Range("E10").Value = "'  +2.98%  "
$ws.Range("E11").Value = "'  +1.03%  "
$ws.Range("E12").Value = "'  +5.02%  "
$ws.Range("D13").Value = "'4.051.13"
$ws.Range("E13").Value = "'  +1.29%  "
$ws.Range("D14").Value = "'29.25"
$ws.Range("E14").Value = "'  +3.11%  "
$ws.Range("E15").Value = "'  +2.46%  "
$ws.Range("D16").Value = "'3.454.18"
$ws.Range("E16").Value = "'  +1.25%  "
$ws.Range("E17").Value = "'  +1.12%  "
$ws.Range("D18").Value = "'63.154.70"
$ws.Range("E18").Value = "'  +2.03%  "
$ws.Range("E19").Value = "'  +3.99%  "
$ws.Range("D20").Value = "'14.44"
$ws.Range("E20").Value = "'  +3.13%  "
$ws.Range("D21").Value = "'9.28"
$ws.Range("E21").Value = "'  +1.35%  "
$ws.Range("D22").Value = "'387.94"
$ws.Range("E22").Value = "'  -0.62%  "
$ws.Range("D23").Value = "'0.562"
$ws.Range("E23").Value = "'  +1.62%  "
$ws.Range("D24").Value = "'74.49"
$ws.Range("E24").Value = "'  -0.63%  "
$ws.Range("E25").Value = "'  +0.03%  "
$ws.Range("D26").Value = "'3.605.04"
$ws.Range("E26").Value = "'  +1.41%  "
$ws.Range("D27").Value = "'0.0000115"
$ws.Range("E27").Value = "'  +0.91%  "
$ws.Range("E28").Value = "'  -1.97%  "
$ws.Range("D29").Value = "'7.63"
$ws.Range("E29").Value = "'  +2.45%  "
$ws.Range("D30").Value = "'0.997"
$ws.Range("E30").Value = "'  -0.16%  "
$ws.Range("D31").Value = "'8.15"
$ws.Range("E31").Value = "'  +1.80%  "
$ws.Range("D32").Value = "'2.12"
$ws.Range("E32").Value = "'  -0.79%  "
$ws.Range("E33").Value = "'  -0.02%  "
$ws.Range("B34").Value = "'Fetch.AI"
$ws.Range("C34").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "'1.34"
$ws.Range("E34").Value = "'  -3.70%  "
$ws.Range("B35").Value = "'EthereumClassic"
$ws.Range("C35").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "'23.39"
$ws.Range("E35").Value = "'  -0.65%  "
$ws.Range("D36").Value = "'7.14"
$ws.Range("E36").Value = "'  +2.50%  "
$ws.Range("E37").Value = "'  +1.42%  "
$ws.Range("B38").Value = "'ImmutableX"
$ws.Range("C38").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'1.59"
$ws.Range("E38").Value = "'  +4.06%  "
$ws.Range("B39").Value = "'EnergySwap"
$ws.Range("C39").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "'31.75"
$ws.Range("E39").Value = "'  +11.27%  "
$ws.Range("D40").Value = "'168.18"
$ws.Range("E40").Value = "'  +0.26%  "
$ws.Range("D41").Value = "'3.494.65"
$ws.Range("E41").Value = "'  +1.42%  "
$ws.Range("D42").Value = "'0.0769"
$ws.Range("E42").Value = "'  +2.06%  "
$ws.Range("D43").Value = "'0.791"
$ws.Range("E43").Value = "'  +0.78%  "
$ws.Range("E45").Value = "'  +3.63%  "
$ws.Range("E46").Value = "'  +3.57%  "
$ws.Range("D47").Value = "'4.36"
$ws.Range("E47").Value = "'  -1.40%  "
$ws.Range("D48").Value = "'2.595.79"
$ws.Range("E48").Value = "'  +3.73%  "
$ws.Range("E49").Value = "'  +11.35%  "
$ws.Range("E50").Value = "'  +2.66%  "
$ws.Range("D51").Value = "'22.95"
